$wb = $excel.ActiveWorkbook

# Update the "Metadata" sheet with the new URL, version, date, and publisher.
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/procedure-modifier"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Update the CodeSystem "System URI" on the "Include from Procedure Modifi" sheet.
$wsCodes = $wb.Worksheets.Item("Include from Procedure Modifi")
$wsCodes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/procedure-modifier"
